$d = $word.ActiveDocument

# --- Change 1: merge "ich nazwy" + bookmark + "," into a single run "ich nazwy," ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$rngFix = $d.Content
[void]$rngFix.Find.Execute("ich nazwy,", $true, $false, $false, $false, $false, $true, 1, $false, "ich nazwy,", 2)

# --- Change 2: replace the two empty paragraphs after "Szczegoly implementacyjne" with 3 new paragraphs ---
$rngHeading = $d.Content
[void]$rngHeading.Find.Execute("Szczegóły implementacyjne", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $rngHeading.Paragraphs(1).Range.End

# Paragraph 1
$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter('Kolejnym krokiem po szczegółowej analizie systemu było zaplanowanie architektury kodu. Naturalnym początkiem było wydzielenie części odpowiedzialnej za przetwarzanie i wyświetlanie nieskompresowanych sekwencji wideo. Po ukończeniu odtwarzacza ruszyły prace, mające na celu stworzenie warstwy widocznej przez użytkownika')
$insertPos = $insertPos + 319

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter(',')
$insertPos = $insertPos + 1

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter(' wraz ')
$insertPos = $insertPos + 6

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter('z funkcją przeprowadzenia wybranych scenariuszy testowych. ')
$insertPos = $insertPos + 59

# Insert paragraph break (consumes one of the pre-existing empty paragraphs)
$rBreak1 = $d.Range($insertPos, $insertPos)
$rBreak1.InsertParagraphAfter()
$insertPos = $insertPos + 1

# Paragraph 2 (before bookmark)
$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter('Biblioteka ')
$insertPos = $insertPos + 11

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter('libVLC')
$italicRange0 = $d.Range($insertPos, $insertPos + 6)
$italicRange0.Font.Italic = $true
$insertPos = $insertPos + 6

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter(' zapewnia dostęp do wielu gotowych metod obsługi materiału wideo. Jednak nie wspiera ona odtwarzania nieskompresowanych sekwencji wideo. Jednym ze standardo')
$insertPos = $insertPos + 156

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter('wych sposobów odtwarzania wideo')
$insertPos = $insertPos + 31

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter(' jest')
$insertPos = $insertPos + 5

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter(',')
$insertPos = $insertPos + 1

# Bookmark _GoBack placed here (zero-length)
$bmRange = $d.Range($insertPos, $insertPos)
[void]$d.Bookmarks.Add("_GoBack", $bmRange)

# Paragraph 2 (after bookmark)
$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter(' poprzez dekompresje filmu, wczytywanie nie całych klatek, a tylko zmian zachodzących między kolejnymi dwoma. Takie podejście pozwala znacząco ograniczyć rozmiar wczytywanych do pamięci danych. ')
$insertPos = $insertPos + 194

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter('Mniejsza ilość danych pozwala na zwiększenie maksymalnej jakości filmu odtwarzanego na tym samym komputerze względem programu ładującego każdą klatkę od nowa. ')
$insertPos = $insertPos + 159

# Insert paragraph break for paragraph 3; then fill the remaining empty paragraph
$pCountBefore = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($pCountBefore)
$insertPos = $lastPara.Range.Start

# Paragraph 3
$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter('W przypadku inkrementacyjnego ładowania zmian w jednym z ostatnich kroków procesu otrzymywana jest klatka reprezentowana w ten sam sposób, co zapisana jako nieskompresowana. Dzięki tej obserwacji, możliwe jest wstrzyknięcie wczytanej nieskompresowanej klatki do standardowego procesu odtwarzania filmu i jego wyświetlenie z użyciem biblioteki ')
$insertPos = $insertPos + 343

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter('libVLC')
$italicRange1 = $d.Range($insertPos, $insertPos + 6)
$italicRange1.Font.Italic = $true
$insertPos = $insertPos + 6

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter('.')
$insertPos = $insertPos + 1

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter(' ')
$insertPos = $insertPos + 1

$rTmp = $d.Range($insertPos, $insertPos)
$rTmp.InsertAfter('Szczegółowy opis implementacji tego sposobu zostanie przedstawiony w dalszej części rozdziału.')
$insertPos = $insertPos + 94

# --- Apply paragraph formatting (firstLine indent 360 twips = 18pt, justify) to the 3 new paragraphs ---
$totalParas = $d.Paragraphs.Count
for ($i = $totalParas - 2; $i -le $totalParas; $i++) {
    $p = $d.Paragraphs($i)
    $p.Format.FirstLineIndent = 18
    $p.Format.Alignment = 3
}
